$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the swapped C/D values on rows 4 and 5 ---
$ws.Range("C4").Value = 20
$ws.Range("C5").Value = 8

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "10.928962"
$ws.Range("D4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "4.371585"
$ws.Range("D5").ClearFormats()

# --- New column E: "Аккум.доля" (accumulated share) values ---
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "54.644809"
$ws.Range("E2").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "81.967213"
$ws.Range("E3").ClearFormats()

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "92.896175"
$ws.Range("E4").ClearFormats()

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "97.267760"
$ws.Range("E5").ClearFormats()

$ws.Range("E6").Value = 150

# --- New column F: "Категория" (category) ---
$ws.Range("F1").Value = "Категория"

$ws.Range("F2").Value = "A"
$ws.Range("F2").Font.Name = "Arial"
$ws.Range("F2").Font.Size = 10

$ws.Range("F3").Value = "B"
$ws.Range("F3").Font.Name = "Arial"
$ws.Range("F3").Font.Size = 10

$ws.Range("F4").Value = "B"
$ws.Range("F4").Font.Name = "Arial"
$ws.Range("F4").Font.Size = 10

$ws.Range("F5").Value = "C"
$ws.Range("F5").Font.Name = "Arial"
$ws.Range("F5").Font.Size = 10

$ws.Range("F6").Value = "C"
$ws.Range("F6").Font.Name = "Arial"
$ws.Range("F6").Font.Size = 10

# --- New (blank, formatted) cell on row 8 ---
$ws.Range("C8").Font.Name = "Arial"
$ws.Range("C8").Font.Size = 10

# --- Column widths for D and E ---
$ws.Columns.Item(4).ColumnWidth = 13.109375
$ws.Columns.Item(5).ColumnWidth = 12.109375

# --- Selection, matching the final state in the workbook ---
$ws.Range("L4").Select() | Out-Null
